$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 42

# Column A holds the date as literal text (matching existing rows), not an
# Excel date value, so force text formatting before assigning the value and
# then clear the formatting so no extra style is introduced.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "01/05/2026"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = 13601.81
$ws.Cells.Item($row, 3).Value = 0.2037872448926116
$ws.Cells.Item($row, 4).Value = 0.7962127551073884
$ws.Cells.Item($row, 5).Value = -85.65000000000001
$ws.Cells.Item($row, 6).Value = -14.64
$ws.Cells.Item($row, 7).Value = -19811.82
$ws.Cells.Item($row, 8).Value = -64.66
$ws.Cells.Item($row, 9).Value = -331.24
$ws.Cells.Item($row, 10).Value = -10.68
